$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.650391101837158
$ws.Range("B1").Value = 3.713689088821411
$ws.Range("C1").Value = 3.281786441802979
$ws.Range("D1").Value = 4.116365432739258
$ws.Range("E1").Value = 5.256838798522949
